$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values look numeric (e.g. "1.001", "30.531.95") but must
# remain literal text, matching the workbook's existing inlineStr cells.
# Force text format, assign, then reset the style so no stray number format
# / style index lingers on the cell (matches original unstyled cells).
function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue 2 4 '30.531.95'
$ws.Cells.Item(2, 5).Value = '  -0.18%  '
Set-TextValue 3 4 '1.917.61'
$ws.Cells.Item(3, 5).Value = '  -0.46%  '
$ws.Cells.Item(4, 5).Value = '  -0.01%  '
Set-TextValue 5 4 '245.48'
$ws.Cells.Item(5, 5).Value = '  -0.34%  '
Set-TextValue 6 4 '1.001'
$ws.Cells.Item(6, 5).Value = '  +0.00%  '
Set-TextValue 7 4 '0.4827'
$ws.Cells.Item(7, 5).Value = '  +1.86%  '
Set-TextValue 8 4 '0.2905'
$ws.Cells.Item(8, 5).Value = '  -0.73%  '
Set-TextValue 9 4 '0.06710'
$ws.Cells.Item(9, 5).Value = '  -1.54%  '
Set-TextValue 10 4 '111.01'
$ws.Cells.Item(10, 5).Value = '  +5.16%  '
Set-TextValue 11 4 '18.92'
$ws.Cells.Item(11, 5).Value = '  +2.69%  '
Set-TextValue 12 4 '1.918.97'
$ws.Cells.Item(12, 5).Value = '  -0.38%  '
$ws.Cells.Item(13, 5).Value = '  -2.45%  '
Set-TextValue 14 4 '5.289'
$ws.Cells.Item(14, 5).Value = '  -1.19%  '
Set-TextValue 15 4 '0.6692'
$ws.Cells.Item(15, 5).Value = '  -0.44%  '
Set-TextValue 16 4 '298.73'
$ws.Cells.Item(16, 5).Value = '  +3.96%  '
Set-TextValue 17 4 '30.539.97'
$ws.Cells.Item(17, 5).Value = '  -0.33%  '
Set-TextValue 18 4 '12.99'
$ws.Cells.Item(18, 5).Value = '  -0.80%  '
Set-TextValue 19 4 '1.001'
$ws.Cells.Item(19, 5).Value = '  +0.02%  '
Set-TextValue 20 4 '0.000007575'
$ws.Cells.Item(20, 5).Value = '  -1.08%  '
Set-TextValue 21 4 '5.553'
$ws.Cells.Item(21, 5).Value = '  +2.18%  '
Set-TextValue 22 4 '2.170.51'
$ws.Cells.Item(22, 5).Value = '  +0.15%  '
Set-TextValue 23 4 '1.002'
$ws.Cells.Item(23, 5).Value = '  +0.00%  '
Set-TextValue 24 4 '6.421'
$ws.Cells.Item(24, 5).Value = '  +2.28%  '
Set-TextValue 25 4 '9.461'
$ws.Cells.Item(25, 5).Value = '  +0.55%  '
Set-TextValue 26 4 '165.25'
$ws.Cells.Item(26, 5).Value = '  -2.01%  '
Set-TextValue 27 4 '20.24'
$ws.Cells.Item(27, 5).Value = '  -2.67%  '
Set-TextValue 28 4 '2.109'
$ws.Cells.Item(28, 5).Value = '  -1.40%  '
Set-TextValue 29 4 '0.1063'
$ws.Cells.Item(29, 5).Value = '  -2.28%  '
Set-TextValue 30 4 '1.431'
$ws.Cells.Item(30, 5).Value = '  +5.19%  '
Set-TextValue 31 4 '4.141'
$ws.Cells.Item(31, 5).Value = '  -0.90%  '
Set-TextValue 32 4 '4.070'
$ws.Cells.Item(33, 5).Value = '  -1.26%  '
Set-TextValue 34 4 '0.7384'
$ws.Cells.Item(34, 5).Value = '  -0.31%  '
Set-TextValue 35 4 '1.136'
$ws.Cells.Item(35, 5).Value = '  -1.56%  '
Set-TextValue 36 4 '1.000'
$ws.Cells.Item(36, 5).Value = '  +0.00%  '
Set-TextValue 37 4 '2.722'
$ws.Cells.Item(37, 5).Value = '  -0.32%  '
Set-TextValue 38 4 '0.02024'
$ws.Cells.Item(38, 5).Value = '  -3.41%  '
$ws.Cells.Item(39, 5).Value = '  -0.45%  '
Set-TextValue 40 4 '110.68'
$ws.Cells.Item(40, 5).Value = '  -0.44%  '
Set-TextValue 41 4 '2.014'
$ws.Cells.Item(41, 5).Value = '  -2.66%  '
Set-TextValue 42 4 '0.4441'
$ws.Cells.Item(42, 5).Value = '  -0.11%  '
Set-TextValue 43 4 '0.8650'
$ws.Cells.Item(43, 5).Value = '  -1.31%  '
Set-TextValue 44 4 '70.97'
$ws.Cells.Item(44, 5).Value = '  +4.80%  '
Set-TextValue 45 4 '5.832'
$ws.Cells.Item(45, 5).Value = '  -1.51%  '
Set-TextValue 46 4 '1.001'
$ws.Cells.Item(46, 5).Value = '  +0.00%  '
Set-TextValue 47 4 '48.90'
$ws.Cells.Item(47, 5).Value = '  +3.40%  '
Set-TextValue 48 4 '7.204'
$ws.Cells.Item(48, 5).Value = '  -1.08%  '
Set-TextValue 49 4 '9.242'
$ws.Cells.Item(49, 5).Value = '  -1.07%  '
Set-TextValue 50 4 '0.1230'
$ws.Cells.Item(50, 5).Value = '  -0.69%  '
Set-TextValue 51 4 '34.84'
$ws.Cells.Item(51, 5).Value = '  -1.32%  '
